$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: GO 2015
$ws.Cells.Item(2, 1).Value = "GO"
$ws.Cells.Item(2, 2).Value = 2015
$ws.Cells.Item(2, 3).Value = "Gau"
$ws.Cells.Item(2, 4).Value = 0.0322
$ws.Cells.Item(2, 5).Value = 1038.8228
$ws.Cells.Item(2, 6).Value = 1236.24
$ws.Cells.Item(2, 7).Value = 0.00003099662425584037
$ws.Cells.Item(2, 8).Value = 0.0232
$ws.Cells.Item(2, 9).Value = -1.2087

# Row 3: GO 2016
$ws.Cells.Item(3, 1).Value = "GO"
$ws.Cells.Item(3, 2).Value = 2016
$ws.Cells.Item(3, 3).Value = "Exp"
$ws.Cells.Item(3, 4).Value = 0.0194
$ws.Cells.Item(3, 5).Value = 0.4549
$ws.Cells.Item(3, 6).Value = 240.19
$ws.Cells.Item(3, 7).Value = 0.0426467355462739
$ws.Cells.Item(3, 8).Value = 0.0061
$ws.Cells.Item(3, 9).Value = 0.2022

# Row 4: GO 2017
$ws.Cells.Item(4, 1).Value = "GO"
$ws.Cells.Item(4, 2).Value = 2017
$ws.Cells.Item(4, 3).Value = "Gau"
$ws.Cells.Item(4, 4).Value = 0.0385
$ws.Cells.Item(4, 5).Value = 0.0575
$ws.Cells.Item(4, 6).Value = 2.66
$ws.Cells.Item(4, 7).Value = 0.6695652173913043
$ws.Cells.Item(4, 8).Value = 0.0067
$ws.Cells.Item(4, 9).Value = -0.2834

# Row 5: GO 2018
$ws.Cells.Item(5, 1).Value = "GO"
$ws.Cells.Item(5, 2).Value = 2018
$ws.Cells.Item(5, 3).Value = "Exp"
$ws.Cells.Item(5, 4).Value = 0.0187
$ws.Cells.Item(5, 5).Value = 0.4214
$ws.Cells.Item(5, 6).Value = 277.87
$ws.Cells.Item(5, 7).Value = 0.04437588989084006
$ws.Cells.Item(5, 8).Value = 0.0033
$ws.Cells.Item(5, 9).Value = -0.3111

# Row 6: GO 2019
$ws.Cells.Item(6, 1).Value = "GO"
$ws.Cells.Item(6, 2).Value = 2019
$ws.Cells.Item(6, 3).Value = "Gau"
$ws.Cells.Item(6, 4).Value = 0.0684
$ws.Cells.Item(6, 5).Value = 0.1702
$ws.Cells.Item(6, 6).Value = 8.49
$ws.Cells.Item(6, 7).Value = 0.4018801410105758
$ws.Cells.Item(6, 8).Value = 0.048
$ws.Cells.Item(6, 9).Value = -0.4541

# Row 7: GO 2020
$ws.Cells.Item(7, 1).Value = "GO"
$ws.Cells.Item(7, 2).Value = 2020
$ws.Cells.Item(7, 3).Value = "Gau"
$ws.Cells.Item(7, 4).Value = 0.1362
$ws.Cells.Item(7, 5).Value = 2190.6024
$ws.Cells.Item(7, 6).Value = 1152.64
$ws.Cells.Item(7, 7).Value = 0.00006217467852678331
$ws.Cells.Item(7, 8).Value = 0.0849
$ws.Cells.Item(7, 9).Value = -11.57

# Row 8: MG 2015
$ws.Cells.Item(8, 1).Value = "MG"
$ws.Cells.Item(8, 2).Value = 2015
$ws.Cells.Item(8, 3).Value = "Exp"
$ws.Cells.Item(8, 4).Value = 0.0269
$ws.Cells.Item(8, 5).Value = 0.575
$ws.Cells.Item(8, 6).Value = 321.68
$ws.Cells.Item(8, 7).Value = 0.04678260869565218
$ws.Cells.Item(8, 8).Value = 0.0323
$ws.Cells.Item(8, 9).Value = -0.1488

# Row 9: MG 2016
$ws.Cells.Item(9, 1).Value = "MG"
$ws.Cells.Item(9, 2).Value = 2016
$ws.Cells.Item(9, 3).Value = "Exp"
$ws.Cells.Item(9, 4).Value = 0.0247
$ws.Cells.Item(9, 5).Value = 0.2812
$ws.Cells.Item(9, 6).Value = 286.22
$ws.Cells.Item(9, 7).Value = 0.08783783783783783
$ws.Cells.Item(9, 8).Value = 0.0037
$ws.Cells.Item(9, 9).Value = -1.1299

# Row 10: MG 2017
$ws.Cells.Item(10, 1).Value = "MG"
$ws.Cells.Item(10, 2).Value = 2017
$ws.Cells.Item(10, 3).Value = "Sph"
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0.0628
$ws.Cells.Item(10, 6).Value = 1.51
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0.0062
$ws.Cells.Item(10, 9).Value = -0.3295

# Row 11: MG 2018
$ws.Cells.Item(11, 1).Value = "MG"
$ws.Cells.Item(11, 2).Value = 2018
$ws.Cells.Item(11, 3).Value = "Exp"
$ws.Cells.Item(11, 4).Value = 0.0181
$ws.Cells.Item(11, 5).Value = 0.3724
$ws.Cells.Item(11, 6).Value = 193.37
$ws.Cells.Item(11, 7).Value = 0.04860365198711063
$ws.Cells.Item(11, 8).Value = 0.0049
$ws.Cells.Item(11, 9).Value = -0.1237

# Row 12: MG 2019
$ws.Cells.Item(12, 1).Value = "MG"
$ws.Cells.Item(12, 2).Value = 2019
$ws.Cells.Item(12, 3).Value = "Gau"
$ws.Cells.Item(12, 4).Value = 0.0225
$ws.Cells.Item(12, 5).Value = 0.0751
$ws.Cells.Item(12, 6).Value = 8.94
$ws.Cells.Item(12, 7).Value = 0.2996005326231691
$ws.Cells.Item(12, 8).Value = 0.0025
$ws.Cells.Item(12, 9).Value = 0.67

# Row 13: MG 2020
$ws.Cells.Item(13, 1).Value = "MG"
$ws.Cells.Item(13, 2).Value = 2020
$ws.Cells.Item(13, 3).Value = "Gau"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0.3547
$ws.Cells.Item(13, 6).Value = 1.26
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0.2104
$ws.Cells.Item(13, 9).Value = -0.0008

# Row 14: MS 2015
$ws.Cells.Item(14, 1).Value = "MS"
$ws.Cells.Item(14, 2).Value = 2015
$ws.Cells.Item(14, 3).Value = "Gau"
$ws.Cells.Item(14, 4).Value = 0.0391
$ws.Cells.Item(14, 5).Value = 0.0408
$ws.Cells.Item(14, 6).Value = 0.87
$ws.Cells.Item(14, 7).Value = 0.9583333333333334
$ws.Cells.Item(14, 8).Value = 0.0003
$ws.Cells.Item(14, 9).Value = -0.2209

# Row 15: MS 2016
$ws.Cells.Item(15, 1).Value = "MS"
$ws.Cells.Item(15, 2).Value = 2016
$ws.Cells.Item(15, 3).Value = "Exp"
$ws.Cells.Item(15, 4).Value = 0.0151
$ws.Cells.Item(15, 5).Value = 0.8643
$ws.Cells.Item(15, 6).Value = 1240.14
$ws.Cells.Item(15, 7).Value = 0.01747078560684947
$ws.Cells.Item(15, 8).Value = 0.0004
$ws.Cells.Item(15, 9).Value = 0.0711

# Row 16: MS 2017
$ws.Cells.Item(16, 1).Value = "MS"
$ws.Cells.Item(16, 2).Value = 2017
$ws.Cells.Item(16, 3).Value = "Sph"
$ws.Cells.Item(16, 4).Value = 0.0285
$ws.Cells.Item(16, 5).Value = 0.1778
$ws.Cells.Item(16, 6).Value = 25.07
$ws.Cells.Item(16, 7).Value = 0.1602924634420697
$ws.Cells.Item(16, 8).Value = 0.0224
$ws.Cells.Item(16, 9).Value = 0.2542

# Row 17: MS 2018
$ws.Cells.Item(17, 1).Value = "MS"
$ws.Cells.Item(17, 2).Value = 2018
$ws.Cells.Item(17, 3).Value = "Exp"
$ws.Cells.Item(17, 4).Value = 0.0142
$ws.Cells.Item(17, 5).Value = 0.4584
$ws.Cells.Item(17, 6).Value = 551.42
$ws.Cells.Item(17, 7).Value = 0.03097731239092496
$ws.Cells.Item(17, 8).Value = 0.0019
$ws.Cells.Item(17, 9).Value = -3.2505

# Row 18: MS 2019
$ws.Cells.Item(18, 1).Value = "MS"
$ws.Cells.Item(18, 2).Value = 2019
$ws.Cells.Item(18, 3).Value = "Exp"
$ws.Cells.Item(18, 4).Value = 0.0126
$ws.Cells.Item(18, 5).Value = 0.075
$ws.Cells.Item(18, 6).Value = 8.13
$ws.Cells.Item(18, 7).Value = 0.168
$ws.Cells.Item(18, 8).Value = 0.0017
$ws.Cells.Item(18, 9).Value = 0.4811

# Row 19: MS 2020
$ws.Cells.Item(19, 1).Value = "MS"
$ws.Cells.Item(19, 2).Value = 2020
$ws.Cells.Item(19, 3).Value = "Exp"
$ws.Cells.Item(19, 4).Value = 0.0132
$ws.Cells.Item(19, 5).Value = 0.0674
$ws.Cells.Item(19, 6).Value = 15.19
$ws.Cells.Item(19, 7).Value = 0.1958456973293768
$ws.Cells.Item(19, 8).Value = 65.9798
$ws.Cells.Item(19, 9).Value = -5.4748

# Row 20: MT 2015
$ws.Cells.Item(20, 1).Value = "MT"
$ws.Cells.Item(20, 2).Value = 2015
$ws.Cells.Item(20, 3).Value = "Exp"
$ws.Cells.Item(20, 4).Value = 0.0219
$ws.Cells.Item(20, 5).Value = 0.4571
$ws.Cells.Item(20, 6).Value = 447.34
$ws.Cells.Item(20, 7).Value = 0.047910741632028
$ws.Cells.Item(20, 8).Value = 0.0048
$ws.Cells.Item(20, 9).Value = -0.7626

# Row 21: MT 2016
$ws.Cells.Item(21, 1).Value = "MT"
$ws.Cells.Item(21, 2).Value = 2016
$ws.Cells.Item(21, 3).Value = "Gau"
$ws.Cells.Item(21, 4).Value = 0.0171
$ws.Cells.Item(21, 5).Value = 0.029
$ws.Cells.Item(21, 6).Value = 2.34
$ws.Cells.Item(21, 7).Value = 0.5896551724137931
$ws.Cells.Item(21, 8).Value = 0.002
$ws.Cells.Item(21, 9).Value = -0.1081

# Row 22: MT 2017
$ws.Cells.Item(22, 1).Value = "MT"
$ws.Cells.Item(22, 2).Value = 2017
$ws.Cells.Item(22, 3).Value = "Exp"
$ws.Cells.Item(22, 4).Value = 0.0209
$ws.Cells.Item(22, 5).Value = 3.9107
$ws.Cells.Item(22, 6).Value = 2730.17
$ws.Cells.Item(22, 7).Value = 0.00534431176004295
$ws.Cells.Item(22, 8).Value = 0.0021
$ws.Cells.Item(22, 9).Value = 0.2872

# Row 23: MT 2018
$ws.Cells.Item(23, 1).Value = "MT"
$ws.Cells.Item(23, 2).Value = 2018
$ws.Cells.Item(23, 3).Value = "Exp"
$ws.Cells.Item(23, 4).Value = 0.0332
$ws.Cells.Item(23, 5).Value = 0.2548
$ws.Cells.Item(23, 6).Value = 350.21
$ws.Cells.Item(23, 7).Value = 0.130298273155416
$ws.Cells.Item(23, 8).Value = 0.0067
$ws.Cells.Item(23, 9).Value = -2.0969

# Row 24: MT 2019
$ws.Cells.Item(24, 1).Value = "MT"
$ws.Cells.Item(24, 2).Value = 2019
$ws.Cells.Item(24, 3).Value = "Gau"
$ws.Cells.Item(24, 4).Value = 0.0393
$ws.Cells.Item(24, 5).Value = 218.8241
$ws.Cells.Item(24, 6).Value = 1124.81
$ws.Cells.Item(24, 7).Value = 0.0001795963058913529
$ws.Cells.Item(24, 8).Value = 0.0066
$ws.Cells.Item(24, 9).Value = -0.0325

# Row 25: MT 2020
$ws.Cells.Item(25, 1).Value = "MT"
$ws.Cells.Item(25, 2).Value = 2020
$ws.Cells.Item(25, 3).Value = "Gau"
$ws.Cells.Item(25, 4).Value = 0.2031
$ws.Cells.Item(25, 5).Value = 698.9262
$ws.Cells.Item(25, 6).Value = 583.16
$ws.Cells.Item(25, 7).Value = 0.0002905886200860692
$ws.Cells.Item(25, 8).Value = 0.0816
$ws.Cells.Item(25, 9).Value = 0.5056

# Row 26: PA 2015
$ws.Cells.Item(26, 1).Value = "PA"
$ws.Cells.Item(26, 2).Value = 2015
$ws.Cells.Item(26, 3).Value = "Sph"
$ws.Cells.Item(26, 4).Value = 0.0181
$ws.Cells.Item(26, 5).Value = 0.0314
$ws.Cells.Item(26, 6).Value = 19.18
$ws.Cells.Item(26, 7).Value = 0.5764331210191084
$ws.Cells.Item(26, 8).Value = 0.0002
$ws.Cells.Item(26, 9).Value = 0.7635

# Row 27: PA 2016
$ws.Cells.Item(27, 1).Value = "PA"
$ws.Cells.Item(27, 2).Value = 2016
$ws.Cells.Item(27, 3).Value = "Sph"
$ws.Cells.Item(27, 4).Value = 0.0224
$ws.Cells.Item(27, 5).Value = 0.0266
$ws.Cells.Item(27, 6).Value = 6.2
$ws.Cells.Item(27, 7).Value = 0.8421052631578948
$ws.Cells.Item(27, 8).Value = 0.0006
$ws.Cells.Item(27, 9).Value = 0.1653

# Row 28: PA 2017
$ws.Cells.Item(28, 1).Value = "PA"
$ws.Cells.Item(28, 2).Value = 2017
$ws.Cells.Item(28, 3).Value = "Exp"
$ws.Cells.Item(28, 4).Value = 0.0338
$ws.Cells.Item(28, 5).Value = 0.126
$ws.Cells.Item(28, 6).Value = 153.54
$ws.Cells.Item(28, 7).Value = 0.2682539682539682
$ws.Cells.Item(28, 8).Value = 0.0009
$ws.Cells.Item(28, 9).Value = 0.7718

# Row 29: PA 2018
$ws.Cells.Item(29, 1).Value = "PA"
$ws.Cells.Item(29, 2).Value = 2018
$ws.Cells.Item(29, 3).Value = "Sph"
$ws.Cells.Item(29, 4).Value = 0.0302
$ws.Cells.Item(29, 5).Value = 0.0398
$ws.Cells.Item(29, 6).Value = 4.79
$ws.Cells.Item(29, 7).Value = 0.7587939698492462
$ws.Cells.Item(29, 8).Value = 0.0015
$ws.Cells.Item(29, 9).Value = 0.3411

# Row 30: PA 2019
$ws.Cells.Item(30, 1).Value = "PA"
$ws.Cells.Item(30, 2).Value = 2019
$ws.Cells.Item(30, 3).Value = "Sph"
$ws.Cells.Item(30, 4).Value = 0.0159
$ws.Cells.Item(30, 5).Value = 0.0505
$ws.Cells.Item(30, 6).Value = 1.26
$ws.Cells.Item(30, 7).Value = 0.3148514851485149
$ws.Cells.Item(30, 8).Value = 0.0097
$ws.Cells.Item(30, 9).Value = -0.146

# Row 31: PA 2020
$ws.Cells.Item(31, 1).Value = "PA"
$ws.Cells.Item(31, 2).Value = 2020
$ws.Cells.Item(31, 3).Value = "Gau"
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0.2465
$ws.Cells.Item(31, 6).Value = 3.22
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0.0687
$ws.Cells.Item(31, 9).Value = 0.0675
